# Translate the ContosoLearn Market Research document body from English to
# French. Each paragraph's existing content is cleared and replaced with one
# or more runs of French text. Every new run gets French (fr-FR) language
# tagging and the Aptos / Aptos / Times New Roman font set used by the
# translated source document.

$d = $word.ActiveDocument
$nbsp = [char]0x00a0

function Set-ParagraphRuns {
    param(
        [int]$ParaIndex,
        [string[]]$Segments
    )

    $para = $d.Paragraphs.Item($ParaIndex).Range
    $start = $para.Start
    $end = $para.End - 1   # exclude the paragraph mark

    # Wipe the existing content of the paragraph (all runs/fields/proofErr).
    if ($end -gt $start) {
        $full = $d.Range($start, $end)
        $full.Text = ""
    }

    # Insert each segment as its own run, then stamp that run's formatting
    # individually so that adjoining runs do not get auto-merged even though
    # their resulting formatting is identical (mirrors the source diff, which
    # keeps several runs per paragraph with matching rPr).
    $cursor = $start
    foreach ($seg in $Segments) {
        $insPoint = $d.Range($cursor, $cursor)
        $insPoint.Text = $seg
        $segEnd = $cursor + $seg.Length
        $runRange = $d.Range($cursor, $segEnd)
        $runRange.LanguageID = "fr-FR"
        $runRange.Font.Name = "Aptos"
        $runRange.Font.NameFarEast = "Aptos"
        $runRange.Font.NameBi = "Times New Roman"
        $cursor = $segEnd
    }
}

Set-ParagraphRuns 1 @(
    "Étude de marché ContosoLearn"
)

Set-ParagraphRuns 2 @(
    "AdatumLearn${nbsp}: AdatumLearn est une plateforme d’apprentissage optimisée par l’IA qui utilise l’intelligence artificielle pour enrichir l’eLearning avec des fonctionnalités qui automatisent toute une variété de tâches.",
    " ",
    "Celle-ci est reconnue pour ses fonctionnalités de création de contenu et sa technologie d’apprentissage adaptatif."
)

Set-ParagraphRuns 3 @(
    "AdventureLearn${nbsp}: AdventureLearn est une autre plateforme d’apprentissage basée sur l’IA qui offre des expériences d’apprentissage personnalisées et des recommandations basées sur des données."
)

Set-ParagraphRuns 4 @(
    "AlpineTraining${nbsp}: AlpineTraining est une plateforme d’apprentissage mobile orientée principalement sur le microapprentissage."
)

Set-ParagraphRuns 5 @(
    "Bellows OnDemand${nbsp}: Bellows OnDemand est une solution d’apprentissage complète qui propose de la création de contenu et de la collaboration sociale."
)

Set-ParagraphRuns 6 @(
    "FabrikamLearning${nbsp}: FabrikamLearning fournit une suite de plateformes d’apprentissage qui répondent à différents besoins en matière d’apprentissage."
)

Set-ParagraphRuns 7 @(
    "FirstUp Cards${nbsp}: FirstUp Cards est une application d’apprentissage mobile qui est idéale pour des formations sur les procédures de sécurité, la conformité, les nouvelles connaissances de produits ou tout autre type de scénario de formation."
)

Set-ParagraphRuns 8 @(
    "Munson’sLearn${nbsp}: Munson’sLearn est conçu pour permettre aux entreprises de former leurs employés, leurs partenaires et leurs clients."
)

Set-ParagraphRuns 9 @(
    "LibertyLearn${nbsp}: LibertyLearn est un système de gestion de formations rapide pour votre projet stratégique."
)

Set-ParagraphRuns 10 @(
    "WoodgroveLMS${nbsp}: WoodgroveLMS est un système de gestion de formations fonctionnel et attrayant conçu pour offrir la meilleure expérience de formation possible."
)

Set-ParagraphRuns 11 @(
    "NorthwindWorlds${nbsp}: NorthwindWorlds est une solution de formation puissante, facile à utiliser et fiable pour les particuliers et les entreprises."
)

Set-ParagraphRuns 12 @(
    "ProsewareLearn${nbsp}: ProsewareLearn est une entreprise de formation en ligne qui offre toute une variété de cours sous forme de formations vidéo pour les développeurs de logiciels, les administrateurs informatiques et les professionnels de la création via son site web."
)

Set-ParagraphRuns 13 @(
    "RelecloudLearn${nbsp}: RelecloudLearn est une plateforme de formation en ligne américaine qui offre des cours en ligne ouverts massivement (MOOC), des spécialisations et des diplômes pour toute une variété de sujets."
)

Set-ParagraphRuns 14 @(
    "TreyAcademy${nbsp}: TreyAcademy est une plateforme de formation en ligne destinée aux adultes professionnels et aux étudiants, développée en mai 2010."
)

Set-ParagraphRuns 15 @(
    "Ces plateformes ont une présence importante sur le marché et sont largement reconnues pour leurs fonctionnalités basées sur l’IA, telles que des expériences d’apprentissage personnalisées, des recommandations basées sur des données et l’automatisation des tâches.",
    " ",
    "Celles-ci transforment le paysage de l’eLearning en tirant parti de l’IA pour offrir des expériences d’apprentissage plus attrayantes, plus enrichissantes et plus personnalisées.",
    "$nbsp"
)

Write-Output "done"
